# Auto-generated edit script: applies per-cell value updates derived from the
# canonical OOXML diff across the 8 profession worksheets (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR). Each row's H:N block (price/profit columns) is refreshed
# with newly-computed market data; some cells are cleared (no longer populated)
# and a few previously-empty cells gain new values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H13").Value = 30000
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 30000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 30000
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -30338
$ws.Range("H53").Value = 500.5
$ws.Range("I53").Value = 489.6111
$ws.Range("J53").Value = 549.5
$ws.Range("K53").Value = 489.6111
$ws.Range("L53").Value = 549.5
$ws.Range("M53").Value = 147.3889
$ws.Range("N53").Value = -1823.5
$ws.Range("H62").Value = 5796.2666
$ws.Range("J62").Value = 5210.2856
$ws.Range("L62").Value = 5210.2856
$ws.Range("N62").Value = -6458.2856
$ws.Range("H65").Value = 5796.2666
$ws.Range("J65").Value = 5210.2856
$ws.Range("L65").Value = 26051.428
$ws.Range("N65").Value = -32291.428
$ws.Range("H98").Value = 1318.7273
$ws.Range("I98").Value = 1400.7
$ws.Range("J98").Value = 499
$ws.Range("K98").Value = 1400.7
$ws.Range("L98").Value = 499
$ws.Range("M98").Value = 97.29999999999995
$ws.Range("N98").Value = -3495
$ws.Range("H122").Value = 1318.7273
$ws.Range("I122").Value = 1400.7
$ws.Range("J122").Value = 499
$ws.Range("K122").Value = 4202.1
$ws.Range("L122").Value = 1497
$ws.Range("M122").Value = -1752.1
$ws.Range("N122").Value = -6397
$ws.Range("H131").Value = 500
$ws.Range("I131").Value = 500
$ws.Range("K131").Value = 1500
$ws.Range("M131").Value = 3540
$ws.Range("H132").Value = 2171.4473
$ws.Range("I132").Value = 1278.129
$ws.Range("K132").Value = 3834.387
$ws.Range("M132").Value = -1304.387

$ws = $wb.Worksheets("ARM")
$ws.Range("H61").Value = 5459.6
$ws.Range("I61").Value = 6324.5
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 6324.5
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -6112.5
$ws.Range("N61").Value = -2424
$ws.Range("H122").Value = 1251.2
$ws.Range("I122").Value = 1251.2
$ws.Range("K122").Value = 3753.6
$ws.Range("M122").Value = -1303.6
$ws.Range("H136").Value = 5459.6
$ws.Range("I136").Value = 6324.5
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 18973.5
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -16423.5
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets("BSM")
$ws.Range("H10").Value = 559
$ws.Range("I10").Value = 559
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 559
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -419
$ws.Range("N10").ClearContents()
$ws.Range("H20").Value = 6738
$ws.Range("I20").Value = 3940.5
$ws.Range("J20").Value = 12333
$ws.Range("K20").Value = 3940.5
$ws.Range("L20").Value = 12333
$ws.Range("M20").Value = -3693.5
$ws.Range("N20").Value = -12827

$ws = $wb.Worksheets("CRP")
$ws.Range("H64").Value = 50000
$ws.Range("J64").Value = 50000
$ws.Range("L64").Value = 50000
$ws.Range("N64").Value = -50496
$ws.Range("H67").Value = 50000
$ws.Range("J67").Value = 50000
$ws.Range("L67").Value = 50000
$ws.Range("N67").Value = -51716
$ws.Range("H99").Value = 7511.6665
$ws.Range("I99").Value = 6883.222
$ws.Range("J99").Value = 8454.333000000001
$ws.Range("K99").Value = 6883.222
$ws.Range("L99").Value = 8454.333000000001
$ws.Range("M99").Value = -5385.222
$ws.Range("N99").Value = -11450.333
$ws.Range("H126").Value = 7511.6665
$ws.Range("I126").Value = 6883.222
$ws.Range("J126").Value = 8454.333000000001
$ws.Range("K126").Value = 20649.666
$ws.Range("L126").Value = 25362.999
$ws.Range("M126").Value = -18179.666
$ws.Range("N126").Value = -30302.999

$ws = $wb.Worksheets("CUL")
$ws.Range("H134").Value = 62505116
$ws.Range("I134").Value = 62505116
$ws.Range("K134").Value = 187515348
$ws.Range("M134").Value = -187510278

$ws = $wb.Worksheets("GSM")
$ws.Range("H70").Value = 14296649
$ws.Range("I70").Value = 20013170
$ws.Range("J70").Value = 5349
$ws.Range("K70").Value = 20013170
$ws.Range("L70").Value = 5349
$ws.Range("M70").Value = -20012900
$ws.Range("N70").Value = -5889
$ws.Range("H73").Value = 14296649
$ws.Range("I73").Value = 20013170
$ws.Range("J73").Value = 5349
$ws.Range("K73").Value = 20013170
$ws.Range("L73").Value = 5349
$ws.Range("M73").Value = -20012234
$ws.Range("N73").Value = -7221
$ws.Range("H122").Value = 167461.17
$ws.Range("I122").Value = 167461.17
$ws.Range("K122").Value = 502383.51
$ws.Range("M122").Value = -499933.51

$ws = $wb.Worksheets("LTW")
$ws.Range("H7").Value = 5000
$ws.Range("I7").Value = 5000
$ws.Range("K7").Value = 5000
$ws.Range("M7").Value = -4888
$ws.Range("H20").Value = 14999
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 14999
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 14999
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -15451
$ws.Range("H29").Value = 10749.5
$ws.Range("I29").Value = 10749.5
$ws.Range("K29").Value = 10749.5
$ws.Range("M29").Value = -10454.5
$ws.Range("H40").Value = 7783
$ws.Range("I40").Value = 5659.4
$ws.Range("J40").Value = 10437.5
$ws.Range("K40").Value = 5659.4
$ws.Range("L40").Value = 10437.5
$ws.Range("M40").Value = -5523.4
$ws.Range("N40").Value = -10709.5
$ws.Range("H76").Value = 26548
$ws.Range("J76").Value = 26548
$ws.Range("L76").Value = 26548
$ws.Range("N76").Value = -27224
$ws.Range("H79").Value = 26548
$ws.Range("J79").Value = 26548
$ws.Range("L79").Value = 26548
$ws.Range("N79").Value = -28888
$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 5000
$ws.Range("K126").Value = 15000
$ws.Range("M126").Value = -12530

$ws = $wb.Worksheets("WVR")
$ws.Range("H4").Value = 2800
$ws.Range("J4").Value = 2800
$ws.Range("L4").Value = 2800
$ws.Range("N4").Value = -3026
$ws.Range("H82").Value = 34250
$ws.Range("J82").Value = 34250
$ws.Range("L82").Value = 34250
$ws.Range("N82").Value = -35016
$ws.Range("H85").Value = 34250
$ws.Range("J85").Value = 34250
$ws.Range("L85").Value = 34250
$ws.Range("N85").Value = -36902
$ws.Range("H117").Value = 42499
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 42499
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 42499
$ws.Range("M117").ClearContents()
$ws.Range("N117").Value = -51677
$ws.Range("H122").Value = 717
$ws.Range("I122").Value = 717
$ws.Range("K122").Value = 2151
$ws.Range("M122").Value = 299
